$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new program record as row 3
$ws.Range("A3").Value = 108434
$ws.Range("B3").Value = "Ingeniería en Logística"
$ws.Range("C3").Value = "Ingeniero(a) en logística"
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = "Acuerdo 002 del 14 de febrero de 2019"
$ws.Range("F3").Value = "Consejo directivo"
$ws.Range("G3").Value = "07 de noviembre del 2019"
$ws.Range("H3").Value = 156
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = "Medellín"

# Widen column F to fit the new "instancia_input" text
$ws.Columns("F").ColumnWidth = 29.6

# Update the active selection to the newly added cell
$ws.Range("K3").Select()
